$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:F25 block (columns B,C,D,E,F for rows 2-25)
$bf = New-Object 'object[,]' 24,5
$bf[0,0] = 1.02
$bf[0,1] = 1.094724718943803
$bf[0,2] = 1.087544245363512
$bf[0,3] = 1.105475243550892
$bf[0,4] = 1.106032582824806
$bf[1,0] = 1.02
$bf[1,1] = 1.096547688467868
$bf[1,2] = 1.08894356324363
$bf[1,3] = 1.107209307454973
$bf[1,4] = 1.107659039332691
$bf[2,0] = 1.02
$bf[2,1] = 1.097724861194594
$bf[2,2] = 1.089846631427033
$bf[2,3] = 1.108329238015807
$bf[2,4] = 1.108709253941864
$bf[3,0] = 1.02
$bf[3,1] = 1.098219181798792
$bf[3,2] = 1.09022572064687
$bf[3,3] = 1.108799560857636
$bf[3,4] = 1.10915024648958
$bf[4,0] = 1.02
$bf[4,1] = 1.09830214781211
$bf[4,2] = 1.090289338810487
$bf[4,3] = 1.108878501415792
$bf[4,4] = 1.109224261062268
$bf[5,0] = 1.02
$bf[5,1] = 1.097731468522577
$bf[5,2] = 1.08985169902474
$bf[5,3] = 1.108335524423957
$bf[5,4] = 1.108715148523755
$bf[6,0] = 1.02
$bf[6,1] = 1.095341307743296
$bf[6,2] = 1.088017651132323
$bf[6,3] = 1.10606172618091
$bf[6,4] = 1.106582716433555
$bf[7,0] = 1.02
$bf[7,1] = 1.091110436007784
$bf[7,2] = 1.084767110148246
$bf[7,3] = 1.102038158242861
$bf[7,4] = 1.102807631147321
$bf[8,0] = 1.02
$bf[8,1] = 1.08827611027611
$bf[8,2] = 1.082586860475934
$bf[8,3] = 1.099343650886663
$bf[8,4] = 1.100278418994608
$bf[9,0] = 1.02
$bf[9,1] = 1.087045354563791
$bf[9,2] = 1.081639508764767
$bf[9,3] = 1.098173842995287
$bf[9,4] = 1.09918011343182
$bf[10,0] = 1.02
$bf[10,1] = 1.086587657737852
$bf[10,2] = 1.081287113541977
$bf[10,3] = 1.097738847692417
$bf[10,4] = 1.09877166764213
$bf[11,0] = 1.02
$bf[11,1] = 1.086685859949438
$bf[11,2] = 1.081362726640259
$bf[11,3] = 1.097832177478866
$bf[11,4] = 1.098859302885806
$bf[12,0] = 1.02
$bf[12,1] = 1.087007532305543
$bf[12,2] = 1.081610390073039
$bf[12,3] = 1.098137895960321
$bf[12,4] = 1.099146361177938
$bf[13,0] = 1.02
$bf[13,1] = 1.087205653266318
$bf[13,2] = 1.081762916181576
$bf[13,3] = 1.098326195588453
$bf[13,4] = 1.09932316240447
$bf[14,0] = 1.02
$bf[14,1] = 1.088357716803182
$bf[14,2] = 1.08264966272946
$bf[14,3] = 1.099421221210297
$bf[14,4] = 1.100351242531281
$bf[15,0] = 1.02
$bf[15,1] = 1.089079433494324
$bf[15,2] = 1.08320500604831
$bf[15,3] = 1.100107269407312
$bf[15,4] = 1.100995279209581
$bf[16,0] = 1.02
$bf[16,1] = 1.089500064628912
$bf[16,2] = 1.083528611896882
$bf[16,3] = 1.100507134715239
$bf[16,4] = 1.101370632926381
$bf[17,0] = 1.02
$bf[17,1] = 1.089643432847976
$bf[17,2] = 1.083638899816806
$bf[17,3] = 1.100643428883466
$bf[17,4] = 1.101498568057292
$bf[18,0] = 1.02
$bf[18,1] = 1.08900203479679
$bf[18,2] = 1.083145455794342
$bf[18,3] = 1.100033693554006
$bf[18,4] = 1.100926211508848
$bf[19,0] = 1.02
$bf[19,1] = 1.086912822843869
$bf[19,2] = 1.081537473450673
$bf[19,3] = 1.098047882768789
$bf[19,4] = 1.099061843204006
$bf[20,0] = 1.02
$bf[20,1] = 1.085596123060697
$bf[20,2] = 1.080523533739967
$bf[20,3] = 1.096796559613965
$bf[20,4] = 1.097886820572068
$bf[21,0] = 1.02
$bf[21,1] = 1.086294433017047
$bf[21,2] = 1.081061324986455
$bf[21,3] = 1.097460177033968
$bf[21,4] = 1.098509994495543
$bf[22,0] = 1.02
$bf[22,1] = 1.089037008968486
$bf[22,2] = 1.083172364968769
$bf[22,3] = 1.100066940226511
$bf[22,4] = 1.100957421162861
$bf[23,0] = 1.02
$bf[23,1] = 1.092206573776598
$bf[23,2] = 1.085609736126201
$bf[23,3] = 1.103080425411213
$bf[23,4] = 1.103785727360037
$ws.Range("B2:F25").Value = $bf

# Update I2:N25 block (columns I,J,K,L,M,N for rows 2-25)
$inBlock = New-Object 'object[,]' 24,6
$inBlock[0,0] = 1.054759433496571
$inBlock[0,1] = 1.099534447332846
$inBlock[0,2] = 1.090199235921284
$inBlock[0,3] = 1.108084829822062
$inBlock[0,4] = 1.108640783222091
$inBlock[0,5] = 1.101095912366422
$inBlock[1,0] = 1.05521519482839
$inBlock[1,1] = 1.101020111055283
$inBlock[1,2] = 1.09141674748309
$inBlock[1,3] = 1.109639621278553
$inBlock[1,4] = 1.110088316736458
$inBlock[1,5] = 1.102583685901753
$inBlock[2,0] = 1.055507641793345
$inBlock[2,1] = 1.101978661950653
$inBlock[2,2] = 1.092201639698915
$inBlock[2,3] = 1.110643073628432
$inBlock[2,4] = 1.111022258407817
$inBlock[2,5] = 1.103543598049343
$inBlock[3,0] = 1.055630001372214
$inBlock[3,1] = 1.102380984476992
$inBlock[3,2] = 1.092530919061478
$inBlock[3,3] = 1.111064313976778
$inBlock[3,4] = 1.111414249884779
$inBlock[3,5] = 1.103946491919815
$inBlock[4,0] = 1.055650511889336
$inBlock[4,1] = 1.10244849828425
$inBlock[4,2] = 1.09258616636848
$inBlock[4,3] = 1.111135006596295
$inBlock[4,4] = 1.111480029907623
$inBlock[4,5] = 1.104014101604424
$inBlock[5,0] = 1.055509279061529
$inBlock[5,1] = 1.101984040350995
$inBlock[5,2] = 1.092206042241297
$inBlock[5,3] = 1.110648704647179
$inBlock[5,4] = 1.111027498707277
$inBlock[5,5] = 1.103548984087629
$inBlock[6,0] = 1.054913972215502
$inBlock[6,1] = 1.100037115472286
$inBlock[6,2] = 1.090611309849217
$inBlock[6,3] = 1.108610824718937
$inBlock[6,4] = 1.109130551584239
$inBlock[6,5] = 1.101599294352277
$inBlock[7,0] = 1.053845933083618
$inBlock[7,1] = 1.096584627378029
$inBlock[7,2] = 1.087778400684124
$inBlock[7,3] = 1.104999361257431
$inBlock[7,4] = 1.105766642271563
$inBlock[7,5] = 1.098141903328921
$inBlock[8,0] = 1.053120852458435
$inBlock[8,1] = 1.094267608535265
$inBlock[8,2] = 1.085873870210424
$inBlock[8,3] = 1.102577215331598
$inBlock[8,4] = 1.103509054426151
$inBlock[8,5] = 1.095821594053634
$inBlock[9,0] = 1.052803730179022
$inBlock[9,1] = 1.093260506884041
$inBlock[9,2] = 1.08504527519938
$inBlock[9,3] = 1.101524791859634
$inBlock[9,4] = 1.102527784483331
$inBlock[9,5] = 1.09481306220255
$inBlock[10,0] = 1.052685457386296
$inBlock[10,1] = 1.092885836799993
$inBlock[10,2] = 1.084736897321757
$inBlock[10,3] = 1.101133316790337
$inBlock[10,4] = 1.102162724781046
$inBlock[10,5] = 1.094437860044007
$inBlock[11,0] = 1.052710849080264
$inBlock[11,1] = 1.092966231628609
$inBlock[11,2] = 1.084803072800641
$inBlock[11,3] = 1.101217315013085
$inBlock[11,4] = 1.102241057445581
$inBlock[11,5] = 1.094518369042501
$inBlock[12,0] = 1.052793963523396
$inBlock[12,1] = 1.093229548609377
$inBlock[12,2] = 1.085019796912321
$inBlock[12,3] = 1.10149244387859
$inBlock[12,4] = 1.102497620289163
$inBlock[12,5] = 1.094782059963585
$inBlock[13,0] = 1.052845109389511
$inBlock[13,1] = 1.093391708589879
$inBlock[13,2] = 1.085153247825187
$inBlock[13,3] = 1.101661885457902
$inBlock[13,4] = 1.10265562087466
$inBlock[13,5] = 1.094944450229859
$inBlock[14,0] = 1.053141831861818
$inBlock[14,1] = 1.094334364922522
$inBlock[14,2] = 1.085928777726755
$inBlock[14,3] = 1.102646983741719
$inBlock[14,4] = 1.103574098562297
$inBlock[14,5] = 1.095888445242618
$inBlock[15,0] = 1.053327109094806
$inBlock[15,1] = 1.094924636385384
$inBlock[15,2] = 1.086414189193991
$inBlock[15,3] = 1.103263930941647
$inBlock[15,4] = 1.104149230017831
$inBlock[15,5] = 1.096479554958663
$inBlock[16,0] = 1.053434873854394
$inBlock[16,1] = 1.095268564405174
$inBlock[16,2] = 1.086696943995225
$inBlock[16,3] = 1.103623437680688
$inBlock[16,4] = 1.10448433650891
$inBlock[16,5] = 1.09682397139569
$inBlock[17,0] = 1.053471567387663
$inBlock[17,1] = 1.095385773092357
$inBlock[17,2] = 1.086793292399155
$inBlock[17,3] = 1.103745961577219
$inBlock[17,4] = 1.104598538854201
$inBlock[17,5] = 1.096941346532653
$inBlock[18,0] = 1.05330726211457
$inBlock[18,1] = 1.09486134394173
$inBlock[18,2] = 1.086362148268406
$inBlock[18,3] = 1.103197774458913
$inBlock[18,4] = 1.104087560978573
$inBlock[18,5] = 1.096416172632479
$inBlock[19,0] = 1.052769501667832
$inBlock[19,1] = 1.093152024659425
$inBlock[19,2] = 1.08495599374307
$inBlock[19,3] = 1.101411440809702
$inBlock[19,4] = 1.102422084835635
$inBlock[19,5] = 1.094704425920731
$inBlock[20,0] = 1.052428613783003
$inBlock[20,1] = 1.092073900451104
$inBlock[20,2] = 1.084068406643124
$inBlock[20,3] = 1.100285065455759
$inBlock[20,4] = 1.101371615603567
$inBlock[20,5] = 1.093624770652371
$inBlock[21,0] = 1.052609589791287
$inBlock[21,1] = 1.0926457622474
$inBlock[21,2] = 1.084539267361274
$inBlock[21,3] = 1.100882490012852
$inBlock[21,4] = 1.101928808311211
$inBlock[21,5] = 1.094197444558011
$inBlock[22,0] = 1.053316231050937
$inBlock[22,1] = 1.094889944207572
$inBlock[22,2] = 1.086385664488973
$inBlock[22,3] = 1.103227668801374
$inBlock[22,4] = 1.104115427668884
$inBlock[22,5] = 1.096444813513979
$inBlock[23,0] = 1.054124328652922
$inBlock[23,1] = 1.097479831565423
$inBlock[23,2] = 1.088513537753031
$inBlock[23,3] = 1.10593551217533
$inBlock[23,4] = 1.106638882558813
$inBlock[23,5] = 1.099038378808942
$ws.Range("I2:N25").Value = $inBlock
